$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-13 (columns A-T), reflecting the addition of the
# "Resolving-Mac" cluster to the NATMI ligand-receptor (Artn-Ret) cross-tabulation and
# the refreshed TPM-based expression statistics for every Sending/Target cluster pair.
$data = @(
    @("FAPs","Artn","Ret","ECs",2,0.6666666666666666,0.7368749999999999,2.210625,0.5924430993326582,0.5924430993326582,3,1,4.328495333333334,12.985486,0.2325244328639614,0.2325244328639614,3.18955999875,28.70603998875,0.1377574956764939,0.1377574956764939),
    @("FAPs","Artn","Ret","FAPs",2,0.6666666666666666,0.7368749999999999,2.210625,0.5924430993326582,0.5924430993326582,3,1,7.361539333333333,22.084618,0.3954579193622197,0.3954579193622196,5.424534296249999,48.82080866625,0.2342863154025978,0.2342863154025978),
    @("FAPs","Artn","Ret","MuSCs",2,0.6666666666666666,0.7368749999999999,2.210625,0.5924430993326582,0.5924430993326582,3,1,6.912382333333333,20.737147,0.3713294477689628,0.3713294477689628,5.093561731875,45.842055586875,0.2199915689097287,0.2199915689097287),
    @("FAPs","Artn","Ret","Resolving-Mac",2,0.6666666666666666,0.7368749999999999,2.210625,0.5924430993326582,0.5924430993326582,1,0.3333333333333333,0.012811,0.038433,0.0006882000048562394,0.0006882000048562392,0.009440105625,0.084960950625,0.0004077193438377809,0.0004077193438377807),
    @("MuSCs","Artn","Ret","ECs",2,0.6666666666666666,0.3936963333333334,1.181089,0.3165295008188679,0.3165295008188679,3,1,4.328495333333334,12.985486,0.2325244328639614,0.2325244328639614,1.704112741583778,15.337014674254,0.07360084266262007,0.07360084266262006),
    @("MuSCs","Artn","Ret","FAPs",2,0.6666666666666666,0.3936963333333334,1.181089,0.3165295008188679,0.3165295008188679,3,1,7.361539333333333,22.084618,0.3954579193622197,0.3954579193622196,2.898211043222445,26.083899389002,0.1251740978105915,0.1251740978105915),
    @("MuSCs","Artn","Ret","MuSCs",2,0.6666666666666666,0.3936963333333334,1.181089,0.3165295008188679,0.3165295008188679,3,1,6.912382333333333,20.737147,0.3713294477689628,0.3713294477689628,2.721379579231445,24.492416213083,0.1175367247416557,0.1175367247416557),
    @("MuSCs","Artn","Ret","Resolving-Mac",2,0.6666666666666666,0.3936963333333334,1.181089,0.3165295008188679,0.3165295008188679,1,0.3333333333333333,0.012811,0.038433,0.0006882000048562394,0.0006882000048562392,0.005043643726333334,0.04539279353700001,0.0002178356040006879,0.0002178356040006879),
    @("Resolving-Mac","Artn","Ret","ECs",1,0.3333333333333333,0.113219,0.339657,0.09102739984847392,0.09102739984847392,3,1,4.328495333333334,12.985486,0.2325244328639614,0.2325244328639614,0.4900679131446667,4.410611218302,0.02116609452484744,0.02116609452484744),
    @("Resolving-Mac","Artn","Ret","FAPs",1,0.3333333333333333,0.113219,0.339657,0.09102739984847392,0.09102739984847392,3,1,7.361539333333333,22.084618,0.3954579193622197,0.3954579193622196,0.8334661217806666,7.501195096026,0.03599750614903033,0.03599750614903033),
    @("Resolving-Mac","Artn","Ret","MuSCs",1,0.3333333333333333,0.113219,0.339657,0.09102739984847392,0.09102739984847392,3,1,6.912382333333333,20.737147,0.3713294477689628,0.3713294477689628,0.7826130153976667,7.043517138578999,0.03380115411757839,0.03380115411757838),
    @("Resolving-Mac","Artn","Ret","Resolving-Mac",1,0.3333333333333333,0.113219,0.339657,0.09102739984847392,0.09102739984847392,1,0.3333333333333333,0.012811,0.038433,0.0006882000048562394,0.0006882000048562392,0.001450448609,0.013054037481,0.0000626450570177706,0.00006264505701777059)
)

$rowIndex = 2
foreach ($dataRow in $data) {
    $colIndex = 1
    foreach ($val in $dataRow) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
